$d = $word.ActiveDocument
$full = $d.Content.WordOpenXML

# 1. Update the "Skills" heading text to "Skills and tools"
$full = $full.Replace("<w:t>Skills</w:t>", "<w:t>Skills and tools</w:t>")

# 2. Replace the old two-column skills table with the new tab-separated paragraph list
$tblStart = $full.IndexOf("<w:tbl>")
$tblEndTag = "</w:tbl>"
$tblEnd = $full.IndexOf($tblEndTag) + $tblEndTag.Length
if ($tblStart -lt 0 -or $tblEnd -lt $tblEndTag.Length) {
    throw "Could not locate the skills table in document XML"
}

$newParas = @'
<w:p><w:pPr><w:spacing w:before="0" w:after="469" w:line="259"/><w:ind w:right="358" w:left="0" w:hanging="370"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">CSS</w:t><w:tab/><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">React</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:after="472" w:line="259"/><w:ind w:right="0" w:left="0" w:hanging="370"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">Tailwind</w:t><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">HTML</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:after="702" w:line="259"/><w:ind w:right="358" w:left="0" w:hanging="370"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve"> SCSS</w:t><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">NextJS</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:after="702" w:line="259"/><w:ind w:right="358" w:left="0" w:hanging="370"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">Nodejs</w:t><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">PHP</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:after="702" w:line="259"/><w:ind w:right="358" w:left="0" w:hanging="370"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">Typescript</w:t><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">React Native</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:after="702" w:line="259"/><w:ind w:right="358" w:left="0" w:hanging="370"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">Solidity</w:t><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">JavaScript </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:after="702" w:line="259"/><w:ind w:right="358" w:left="0" w:hanging="370"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">Bootstrap</w:t><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Flutter</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:after="702" w:line="259"/><w:ind w:right="358" w:left="0" w:hanging="370"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato" w:cs="Lato" w:eastAsia="Lato"/><w:b/><w:color w:val="000000"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">Github </w:t><w:tab/><w:tab/><w:tab/><w:t xml:space="preserve">Git</w:t></w:r></w:p>
'@

$full = $full.Substring(0, $tblStart) + $newParas + $full.Substring($tblEnd)

# 3. Renumber the w:numId references (both in the document body numPr and in
#    word/numbering.xml w:num definitions). Go through a unique placeholder per
#    pair first so overlapping old/new values (45/48) do not collide.
$numIdMap = [ordered]@{
    16 = 13
    21 = 18
    26 = 23
    31 = 28
    35 = 32
    39 = 36
    45 = 42
    48 = 45
    51 = 48
    55 = 52
}

$idx = 0
$placeholders = @{}
foreach ($old in $numIdMap.Keys) {
    $placeholder = "@@NUMID_PLACEHOLDER_" + $idx + "@@"
    $placeholders[$old] = $placeholder
    $full = $full.Replace('w:numId w:val="' + $old + '"', 'w:numId w:val="' + $placeholder + '"')
    $full = $full.Replace('<w:num w:numId="' + $old + '">', '<w:num w:numId="' + $placeholder + '">')
    $idx = $idx + 1
}
foreach ($old in $numIdMap.Keys) {
    $new = $numIdMap[$old]
    $placeholder = $placeholders[$old]
    $full = $full.Replace('w:numId w:val="' + $placeholder + '"', 'w:numId w:val="' + $new + '"')
    $full = $full.Replace('<w:num w:numId="' + $placeholder + '">', '<w:num w:numId="' + $new + '">')
}

# 4. Write the whole modified package back into the document.
$d.Content.InsertXML($full)
